$wb = $excel.ActiveWorkbook

$wsSurvey = $wb.Worksheets.Item("survey")

# Row 2: OS_FOL_date type changes from "date" to "text"
$wsSurvey.Range("C2").Value = "text"

# Row 4: OS_time_begin type changes from "time" to "text"
$wsSurvey.Range("C4").Value = "text"

# Row 5: OS_time_end type changes from "time" to "text"
$wsSurvey.Range("C5").Value = "text"

# Update the active selection to C7 (as recorded in the saved view state)
$wsSurvey.Range("C7").Select()
